# Append a new data row (row 67) to each of the four worksheets,
# mirroring the existing data-entry pattern used for the prior rows.

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = "FE_LFT_#1"
        A = [double]"45853.49097222222"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x40"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 320
        I = 15
    },
    @{
        Sheet = "FE_LFT_#2"
        A = [double]"45853.49097222222"
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x50"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 336
        I = 14
    },
    @{
        Sheet = "FE_PLT_#1"
        A = [double]"45853.49097222222"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x64"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 100
        I = 3
    },
    @{
        Sheet = "FE_PLT_#2"
        A = [double]"45853.49097222222"
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x64"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 100
        I = 3
    }
)

foreach ($rowData in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)
    $newRow = 67

    # Column A carries the same date/time number format as the row above it.
    $ws.Range("A$newRow").NumberFormat = $ws.Range("A66").NumberFormat
    $ws.Range("A$newRow").Value = $rowData.A

    $ws.Range("B$newRow").Value = $rowData.B
    $ws.Range("C$newRow").Value = $rowData.C
    $ws.Range("D$newRow").Value = $rowData.D
    $ws.Range("E$newRow").Value = $rowData.E
    $ws.Range("F$newRow").Value = $rowData.F
    $ws.Range("G$newRow").Value = $rowData.G
    $ws.Range("H$newRow").Value = $rowData.H
    $ws.Range("I$newRow").Value = $rowData.I
}
